$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new measurement row for 22/01/2018
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "1/22/2018"

# Column C is formatted as Text, so typing "15.11" stores it as text
# (and trips Excel's "number stored as text" warning).
$ws.Range("C6").Value = "15.11"

$ws.Range("D6").Value = 100.2
$ws.Range("E6").Value = 220.9
$ws.Range("F6").Formula = "=E6-E5"
$ws.Range("G6").Formula = "=E6-210"

$ws.Range("C10").Select()
